# RPA datasets push 2024-05-03
# Insert two new IPO records (디앤디파마텍, 유안타제16호스팩) at the top of the
# data table (rows 2-3), pushing the previously existing 9 rows down by two
# rows (old row 2 -> row 4, ... old row 10 -> row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the existing 9 data rows (2..10) down to (4..12).
# Work from the bottom up so a source row is never overwritten before it is
# copied. Range.Copy(Destination) duplicates the raw cell content/type
# (shared-string index, number, etc.) without re-interpreting text, so the
# already-parsed date-like strings keep being stored as plain text.
for ($r = 10; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 2
    $src = $ws.Range("A" + $srcRow + ":T" + $srcRow)
    $dst = $ws.Range("A" + $dstRow + ":T" + $dstRow)
    $src.Copy($dst)
}

# 2) Populate the two new rows with the new IPO data.
# Columns A (청약일/subscription date), D (납입일/payment date) and
# E (상장일/listing date) hold date-looking text ("2024-04-22" etc.). A plain
# .Value assignment would make Excel auto-convert such text to a real date
# serial number, which also changes the cell's type/style. Formatting the
# cell as Text ("@") first keeps the literal string intact.
$dateCols = @(1, 4, 5)

$newRows = @(
    @{ Row = 2; Values = @("2024-04-22", "디앤디파마텍", "한국", "2024-04-25", "2024-05-02", 36300000, 1100000, "-", 22000, 26000, "-", 33000, "-", "-", 0, "-", "-", "1544 : 1", "-", "-") },
    @{ Row = 3; Values = @("2024-04-22", "유안타제16호스팩", "유안타", "2024-04-25", "2024-05-02", 10300000, 5150000, "-", 2000, 2000, "-", 2000, "-", "-", 0, "-", "-", "334 : 1", "-", "-") }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $values = $entry.Values
    for ($c = 1; $c -le 20; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($dateCols -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $values[$c - 1]
    }
}

# 3) Make sure the sheet dimension covers the two extra rows.
$ws.Range("A1:T12").Select() | Out-Null

Write-Output "RPA datasets push 2024-05-03 applied"
